$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "id" column (A) - data now starts with "name"
$ws.Columns.Item(1).Delete() | Out-Null

# Remove the trailing duplicate "created date" columns (old J:K, now I:J after the shift above)
$ws.Range("I1:J3").Delete() | Out-Null

# Leave the selection where the user last clicked: the empty column right after the data
$ws.Columns.Item(9).Select() | Out-Null
